# "sql alchemy search updated"
# The sheet originally held 9 guest rows (rows 2-10); it is trimmed down to a
# single guest row, with a new name and phone number for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused guest rows (rows 3 through 10), keeping only row 2.
$ws.Range("A3:E10").EntireRow.Delete()

# Replace the remaining guest's name and phone number.
$ws.Range("A2").Value = "rocktim"
$ws.Range("B2").Value = 1423525334

# Match the saved selection state (single cell B2 selected).
$ws.Range("B2").Select()
